$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Pluto"
$ws.Range("B4").Value = "Hyderabad"

$ws.Range("B5").Select()
